# ---------------------------------------------------------------------------
# Edit script: apply the changes described in the commit
#   "added some important risk factors"
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument
$WNS = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-PkgXml([string]$bodyXml) {
    return @"
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document $WNS>
<w:body>
$bodyXml
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
}

# ---------------------------------------------------------------------------
# 1. Aims paragraph: "U.S." -> "New Mexico" (two spots)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(" from the U.S. between", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " from the New Mexico between", 2) | Out-Null

$d.Content.Find.Execute("efficacy of midwife services in the U.S. ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "efficacy of midwife services in the New Mexico. ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Introduction paragraph: same text, but re-written with proofErr markers
#    around "cia", "factbook", "ny", "is" (spelling / grammar check marks).
# ---------------------------------------------------------------------------
$introPara = $d.Paragraphs(7)
$introText = $introPara.Range.Text
if ($introText.StartsWith("One of the main goals of the Affordable Care Act")) {
    $introRange = $d.Range($introPara.Range.Start, $introPara.Range.End - 1)
    $introXml = @"
<w:p>
  <w:r><w:t xml:space="preserve">One of the main goals of the Affordable Care Act was to lower healthcare costs for Americans by encouraging individuals to get covered by health insurance, expanding Medicaid, and encouraging hospitals to lower healthcare cost while maintaining healthcare quality (citation). Labor and delivery should be foremost in hospitals and policy-makers&#8217; minds. Not only does the U.S. have one of the highest infant mortality rates in the developed world (</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>cia</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>factbook</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>), but the cost of child birth in the U.S. is much higher than other Western countries as well (</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>ny</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> times). One significant difference between giving birth in the U.S. and other developed countries </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>is</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t xml:space="preserve"> the percentage of births that are attended by a midwife. In the U.S. only eight percent of deliveries are attended by a midwife as compared to 45 and 68 percent in the Netherlands and Britain, respectively. Although increasing the number of births attended by a midwife may indeed lower costs, because the child-bearing population in the U.S. may look much different than the populations in Europe and other parts of the world, it is important to investigate how birth outcomes in the U.S. vary based on birth attendant.</w:t></w:r>
</w:p>
"@
    $introRange.InsertXML((New-PkgXml $introXml))
}

# ---------------------------------------------------------------------------
# 3. Literature review: re-write the Berglund paragraph (add proofErr marks,
#    same text, no bookmark) and insert three new paragraphs after it
#    (Black/Mitchell citation, a blank paragraph, a paragraph that now only
#    holds the _GoBack bookmark, and the Jena/Prasad citation).
# ---------------------------------------------------------------------------
$bergPara = $d.Paragraphs(9)
$bergText = $bergPara.Range.Text
if ($bergText.StartsWith("Berglund, Lindberg, Nystrom")) {

    # Remove the _GoBack bookmark currently sitting inside this paragraph;
    # it will be re-created later in its own paragraph.
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks.Item("_GoBack").Delete()
    }

    $bergRange = $d.Range($bergPara.Range.Start, $bergPara.Range.End - 1)
    $bergXml = @"
<w:p>
  <w:r><w:t xml:space="preserve">Berglund, Lindberg, Nystrom and </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>LindMark</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> show that there </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>is</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t xml:space="preserve"> no adverse effects on risk assessment when Swedish women&#8217;s risk level is assessed only by a midwife as compared to a midwife and doctor in sequence (2007).</w:t></w:r>
</w:p>
"@
    $bergRange.InsertXML((New-PkgXml $bergXml))

    # Insert the 4 new paragraphs right after the (now bookmark-free) Berglund
    # paragraph: Black/Mitchell citation, blank line, bookmark-only line,
    # Jena/Prasad citation.
    $bergPara2 = $d.Paragraphs(9)
    $insertPos = $bergPara2.Range.End - 1
    $pPrXml = '<w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:firstLine="720"/></w:pPr>'
    $newParasXml = @"
<w:p>
  $pPrXml
  <w:r><w:t xml:space="preserve">Black, Mitchell, and </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>Danielian</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> demonstrate that for one training hospital in the U.K. midwifes have no worse outcomes than </w:t></w:r>
  <w:r><w:t xml:space="preserve">doctors when perform instrument assisted </w:t></w:r>
  <w:r><w:t xml:space="preserve">births. </w:t></w:r>
</w:p>
<w:p>
  $pPrXml
</w:p>
<w:p>
  $pPrXml
</w:p>
<w:p>
  $pPrXml
  <w:r><w:t xml:space="preserve">Jena, Prasad, Goldman and </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>Romley</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> demonstrate that patients treated for AMI and heart-failure at major teaching hospitals during two national conferences for Cardiologists had mortality outcomes no worse than those treated during non-conference periods. </w:t></w:r>
</w:p>
"@
    $d.Range($insertPos, $insertPos).InsertXML((New-PkgXml $newParasXml))

    # Re-create the _GoBack bookmark as its own (now 3rd new) empty paragraph.
    $bookmarkPara = $d.Paragraphs(12)
    $d.Bookmarks.Add("_GoBack", $bookmarkPara.Range) | Out-Null
}

# ---------------------------------------------------------------------------
# 4. References: Tracy et al. paragraph gets the same proofErr treatment.
# ---------------------------------------------------------------------------
$i = 1
$tracyIndex = -1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Tracy, S. K.")) {
        $tracyIndex = $i
    }
    $i = $i + 1
}

if ($tracyIndex -gt 0) {
    $tracyPara = $d.Paragraphs($tracyIndex)
    $tracyRange = $d.Range($tracyPara.Range.Start, $tracyPara.Range.End - 1)
    $tracyXml = @"
<w:p>
  <w:r><w:t xml:space="preserve">Tracy, S. K., </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>Hartz</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve">, D. L., Tracy, M. B., Allen, J., </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>Forti</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve">, A., Hall, B., . . . </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>Kildea</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve">, S. (2013). Caseload midwifery care versus standard maternity care for women of any risk: M@NGO, a </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>randomised</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> controlled trial.</w:t></w:r>
  <w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve"> The Lancet, 382</w:t></w:r>
  <w:r><w:t xml:space="preserve">(9906), 1723-32. </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>doi:http</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>://dx.doi.org/10.1016/S0140-6736(13)61406-3</w:t></w:r>
</w:p>
"@
    $tracyRange.InsertXML((New-PkgXml $tracyXml))
}

Write-Host "Done."
